$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for "07-04-2025" Status/Time columns (Z/AA), copying the
# formatting (bold, border, centered) from the existing header cell X1
# ("01-04-2025 Status") so the new headers match the rest of the header row.
$ws.Range("X1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Z1").Value = "07-04-2025 Status"

$ws.Range("Y1").Copy()
$ws.Range("AA1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AA1").Value = "07-04-2025 Time"

# Fill in attendance data for each student row (rows 2-6): Status "A" and Time "00:00:00"
for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 26).Value = "A"         # column Z = 26
    $ws.Cells.Item($row, 27).Value = "00:00:00"  # column AA = 27
}
